$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.917.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "'2.968.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'592.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "'141.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").Value = "'2.967.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "'5.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.53%  "
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "'0.0000225"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "'33.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "'3.456.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "'61.002.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").Value = "'2.964.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "'448.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'0.676"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'7.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").Value = "'82.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").Value = "'2.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").Value = "'10.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").Value = "'11.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'2.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'7.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").Value = "'2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "'27.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").Value = "'0.0₃0799"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").Value = "'5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").Value = "'50.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").Value = "'8.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("E41").Value = "  +5.42%  "
$ws.Range("D42").Value = "'2.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").Value = "'386.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").Value = "'38.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").Value = "'0.0346"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "'2.686.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'129.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'0.107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").Value = "'2.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
